$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the prior row down to the new row first, so the
# new cells inherit the existing date/left-align styles instead of Excel
# auto-creating a new number format when a date value is assigned.
$ws.Range("A22:C22").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)  # xlPasteFormats

# Append a new time log entry for 11/17 (8 hours, daily operations task)
$ws.Range("A23").Value = 45247
$ws.Range("B23").Value = "Internship"
$ws.Range("C23").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Move selection to C24 (matches Excel's post-edit selection state)
$ws.Range("C24").Select()
